# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated counts.

$wb = $excel.ActiveWorkbook

# Row number -> new value for column F
$updates = @{
    2  = 62
    3  = 1043
    4  = 39
    5  = 72
    6  = 2814
    8  = 1823
    10 = 87
    11 = 648
    12 = 31
    13 = 23
    14 = 180
    15 = 87
    17 = 25
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
